$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing "+TL" techniques (drop the "_optimized" qualifier before "+TL")
$ws.Range("A5").Value = "TFKeras+TL"
$ws.Range("A6").Value = "TFKeras+TL"
$ws.Range("A9").Value = "pytorch+TL"
$ws.Range("A10").Value = "pytorch+TL"
$ws.Range("A15").Value = "pytorch+TL+DLv2"
$ws.Range("A16").Value = "pytorch+TL+DLv2"

# Fix training_time typo on row 16 (1.34 -> 134)
$ws.Range("C16").Value = 134

# Add new row 17 - fine tuning transfer learning technique
$ws.Range("A17").Value = "py_torch+TL_Optimized+DLv2"
$ws.Range("B17").Value = 0.82
$ws.Range("C17").Value = 343
$ws.Range("D17").Value = 0

# Resize the autofilter to cover the new row. Toggling off then on again
# forces the plain-range AutoFilter to re-anchor on the new extent.
$ws.Range("A1:D17").AutoFilter()
$ws.Range("A1:D17").AutoFilter()

# Keep the _FilterDatabase defined name in sync with the autofilter range
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$D`$17"
    }
}

# Update selection to match the saved cursor position
$ws.Range("H14").Select()
